$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.718.49'
$ws.Range("E2").Value = '  -3.44%  '
$ws.Range("D3").Value = '2.924.02'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'548.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("D6").Value = "'130.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("D8").Value = "'0.508"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.64%  '
$ws.Range("D9").Value = '2.918.44'
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("E10").Value = '  -2.87%  '
$ws.Range("E11").Value = '  -3.99%  '
$ws.Range("D12").Value = "'0.444"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.95%  '
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").Value = "'32.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").Value = '3.405.28'
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("D17").Value = "'6.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.18%  '
$ws.Range("D18").Value = '2.918.16'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = '57.652.25'
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("D20").Value = "'416.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.11%  '
$ws.Range("D21").Value = "'13.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.48%  '
$ws.Range("D22").Value = "'0.690"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("D23").Value = "'13.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.00%  '
$ws.Range("D24").Value = "'6.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.59%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = "'0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("D28").Value = "'2.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("E29").Value = '  +3.51%  '
$ws.Range("D30").Value = "'7.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.00%  '
$ws.Range("D31").Value = "'25.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").Value = "'5.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.86%  '
$ws.Range("D35").Value = "'0.929"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("E36").Value = '  +4.42%  '
$ws.Range("D37").Value = "'48.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.00%  '
$ws.Range("E38").Value = '  +3.45%  '
$ws.Range("E39").Value = '  +5.75%  '
$ws.Range("D40").Value = "'2.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.54%  '
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("D42").Value = '2.702.09'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = "'371.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = "'123.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.00%  '
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("D50").Value = "'22.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("E51").Value = '  -0.24%  '
